# Updates cryptos list (coin prices / 1h volume %) per the Tue Aug 29 00:26:39 UTC 2023
# GitHub Actions refresh. Rows 12/13 and 49/50 also swap their ranking order
# (Coin name + Link move along with the row), so B/C/D/E are all rewritten there.
#
# NOTE: Price strings such as "20.67" or "0.07700" are stored as literal TEXT in
# this sheet (not numbers), and some rely on exact formatting (trailing zeros like
# 0.07700, or 3.500) that a numeric cell would collapse. A bare .Value = "20.67"
# would let Excel auto-sense it as a Number and mangle/round it, so for every D-column
# price that looks numeric we prefix it with a literal leading apostrophe (the classic
# "force text" entry method), exactly like typing '20.67 into the cell by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin price/volume refresh
$ws.Range("D2").Value = '26.225.15'
$ws.Range("E2").Value = '  -0.15%  '

# Row 3: Ethereum price/volume refresh
$ws.Range("D3").Value = '1.660.23'
$ws.Range("E3").Value = '  -0.42%  '

# Row 4: TetherUSD price/volume refresh
$ws.Range("E4").Value = '  -0.40%  '

# Row 5: BNB price/volume refresh
$ws.Range("D5").Value = '''219.65'
$ws.Range("E5").Value = '  -0.17%  '

# Row 6: XRP price/volume refresh
$ws.Range("D6").Value = '''0.5258'
$ws.Range("E6").Value = '  -0.37%  '

# Row 8: Cardano price/volume refresh
$ws.Range("D8").Value = '''0.2692'
$ws.Range("E8").Value = '  +1.65%  '

# Row 9: Dogecoin price/volume refresh
$ws.Range("D9").Value = '''0.06388'
$ws.Range("E9").Value = '  +0.49%  '

# Row 10: Solana price/volume refresh
$ws.Range("D10").Value = '''20.67'
$ws.Range("E10").Value = '  -0.84%  '

# Row 11: TRON price/volume refresh
$ws.Range("D11").Value = '''0.07700'
$ws.Range("E11").Value = '  -1.78%  '

# Row 12: ranking swap, now WrappedEther (was Polkadot)
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.887.95'
$ws.Range("E12").Value = '  +12.53%  '

# Row 13: ranking swap, now Polkadot (was WrappedEther)
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.631'
$ws.Range("E13").Value = '  +2.81%  '

# Row 14: WrappedliquidstakedEther2.0 price/volume refresh
$ws.Range("D14").Value = '1.888.92'
$ws.Range("E14").Value = '  -0.44%  '

# Row 15: Polygon price/volume refresh
$ws.Range("D15").Value = '''0.5653'
$ws.Range("E15").Value = '  +1.07%  '

# Row 16: ShibaInu price/volume refresh
$ws.Range("D16").Value = '0.0₅8262'

# Row 17: Litecoin price/volume refresh
$ws.Range("D17").Value = '''65.71'
$ws.Range("E17").Value = '  -0.05%  '

# Row 18: WrappedBTC price/volume refresh
$ws.Range("D18").Value = '26.212.91'
$ws.Range("E18").Value = '  -0.36%  '

# Row 19: Dai price/volume refresh
$ws.Range("E19").Value = '  -0.35%  '

# Row 20: Uniswap price/volume refresh
$ws.Range("D20").Value = '''4.695'
$ws.Range("E20").Value = '  -0.12%  '

# Row 21: Avalanche price/volume refresh
$ws.Range("D21").Value = '''10.47'
$ws.Range("E21").Value = '  +2.17%  '

# Row 22: BitcoinCash price/volume refresh
$ws.Range("D22").Value = '''192.57'
$ws.Range("E22").Value = '  -3.15%  '

# Row 23: Chainlink price/volume refresh
$ws.Range("D23").Value = '''6.005'
$ws.Range("E23").Value = '  -0.60%  '

# Row 24: BinanceUSD price/volume refresh
$ws.Range("E24").Value = '  -0.45%  '

# Row 25: Monero price/volume refresh
$ws.Range("D25").Value = '''145.73'
$ws.Range("E25").Value = '  -0.26%  '

# Row 26: Stellar price/volume refresh
$ws.Range("D26").Value = '''0.1202'
$ws.Range("E26").Value = '  -1.30%  '

# Row 27: Cosmos price/volume refresh
$ws.Range("D27").Value = '''7.310'
$ws.Range("E27").Value = '  +1.20%  '

# Row 28: EthereumClassic price/volume refresh
$ws.Range("D28").Value = '''16.08'
$ws.Range("E28").Value = '  -0.54%  '

# Row 29: Toncoin price/volume refresh
$ws.Range("D29").Value = '''1.525'
$ws.Range("E29").Value = '  -0.07%  '

# Row 30: Hedera price/volume refresh
$ws.Range("D30").Value = '''0.05633'
$ws.Range("E30").Value = '  -4.55%  '

# Row 31: PancakeSwap price/volume refresh
$ws.Range("E31").Value = '  -0.30%  '

# Row 32: InternetComputer(DFINITY) price/volume refresh
$ws.Range("D32").Value = '''3.500'
$ws.Range("E32").Value = '  +0.03%  '

# Row 33: Filecoin price/volume refresh
$ws.Range("D33").Value = '''3.415'
$ws.Range("E33").Value = '  +2.40%  '

# Row 34: LidoDAOToken price/volume refresh
$ws.Range("D34").Value = '''1.582'
$ws.Range("E34").Value = '  -0.82%  '

# Row 35: ARBITRUM price/volume refresh
$ws.Range("D35").Value = '''0.9553'
$ws.Range("E35").Value = '  -1.02%  '

# Row 36: MXToken price/volume refresh
$ws.Range("D36").Value = '''2.786'
$ws.Range("E36").Value = '  -1.28%  '

# Row 37: HuobiToken price/volume refresh
$ws.Range("D37").Value = '''2.404'
$ws.Range("E37").Value = '  -1.02%  '

# Row 38: ImmutableX price/volume refresh
$ws.Range("D38").Value = '''0.5761'
$ws.Range("E38").Value = '  -0.74%  '

# Row 39: VeChain price/volume refresh
$ws.Range("D39").Value = '''0.01603'
$ws.Range("E39").Value = '  -0.17%  '

# Row 40: FraxShare price/volume refresh
$ws.Range("D40").Value = '''5.971'
$ws.Range("E40").Value = '  +0.34%  '

# Row 41: PaxDollar price/volume refresh
$ws.Range("E41").Value = '  -0.37%  '

# Row 42: TrustWalletToken price/volume refresh
$ws.Range("D42").Value = '''0.8387'
$ws.Range("E42").Value = '  -2.48%  '

# Row 43: Maker price/volume refresh
$ws.Range("D43").Value = '1.029.79'
$ws.Range("E43").Value = '  -4.65%  '

# Row 44: Quant price/volume refresh
$ws.Range("D44").Value = '''101.46'
$ws.Range("E44").Value = '  -1.87%  '

# Row 45: RocketPoolETH price/volume refresh
$ws.Range("D45").Value = '1.799.32'
$ws.Range("E45").Value = '  -0.47%  '

# Row 46: Aave price/volume refresh
$ws.Range("D46").Value = '''58.55'
$ws.Range("E46").Value = '  -0.02%  '

# Row 47: BabyDogeCoin price/volume refresh
$ws.Range("D47").Value = '0.0₈106'
$ws.Range("E47").Value = '  -0.93%  '

# Row 48: Frax price/volume refresh
$ws.Range("D48").Value = '''0.9991'
$ws.Range("E48").Value = '  -1.26%  '

# Row 49: ranking swap, now EnergySwap (was Cronos)
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''8.087'
$ws.Range("E49").Value = '  +0.18%  '

# Row 50: ranking swap, now Cronos (was EnergySwap)
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.05306'
$ws.Range("E50").Value = '  +3.15%  '

# Row 51: Mantle price/volume refresh
$ws.Range("E51").Value = '  -1.57%  '
